# SubCalc contour updates, small FIELDS updates
$wb = $excel.ActiveWorkbook
$wsInfo = $wb.Worksheets.Item(1)      # "INFO" sheet
$wsTemplate = $wb.Worksheets.Item(2)  # "SubCalc_template" sheet

# ---------------------------------------------------------------------
# 1) INFO sheet: rewrite the blurb + field descriptions, and insert a
#    new "Group:" row (pushing the rest of the legend down by one row).
# ---------------------------------------------------------------------

$wsInfo.Range("A1").Value = "The SubCalc template should be filled out as a flat file with entries in all columns for all rows, then saved as a csv."

$wsInfo.Range("A3").Value = "Group:"
$wsInfo.Range("B3").Value = "Label identifying groups of footprints that are plotted under the same label. This label will be plotted."

# Row 4 becomes a second bold/formatted legend row (same look as row 3)
# holding the old "Name:" entry.
$wsInfo.Range("A3:B3").Copy()
$wsInfo.Range("A4:B4").PasteSpecial(-4122)   # xlPasteFormats
$wsInfo.Range("A4").Value = "Name:"
$wsInfo.Range("B4").Value = "Label identifiying unique footprints. This label won't be plotted."

$wsInfo.Range("A5").Value = "X:"
$wsInfo.Range("B5").Value = "x coordinates of the footprint. They should be sorted according to the path of the footprint, no jumbling"

$wsInfo.Range("A6").Value = "Y:"
$wsInfo.Range("B6").Value = "y coordinates of the footprint. They should be sorted according to the path of the footprint, no jumbling"

$wsInfo.Range("A7").Value = "Power Line?"
$wsInfo.Range("B7").Value = "A flag indicating whether the footprint corresponds to the modeled power lines (1 = yes, 0 = no)"

$wsInfo.Range("A8").Value = "Of Concern?"
$wsInfo.Range("B8").Value = "A flag indicating whether the structure is of concern with respect to EMF. Things like houses would qualify (1 = yes, 0 = no)"

$wsInfo.Range("A9").Value = "Draw as Loop?"
$wsInfo.Range("B9").Value = "A flag indicating whether the footprint coordinates should be closed upon plotting (1 = yes, 0 = no)"

$wsInfo.Range("F16").Select()

# ---------------------------------------------------------------------
# 2) SubCalc_template sheet: drop the 41 blank template rows that used
#    to pre-fill the flat-file area, leaving just the header row.
# ---------------------------------------------------------------------

$wsTemplate.Range("A2:A42").EntireRow.Delete()
$wsTemplate.Range("F15").Select()
